# Weekly update: a new date-pair (Primera/Segunda) record for
# "Vega Monumental Concepción - Acelga" is inserted at the top of the
# repeating data block (row 43), pushing the existing records down by two
# rows. The new pair duplicates the values that used to sit at rows 43-44,
# except for the reporting date (column D), which becomes 2021-12-30
# (Excel serial 44525).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 43-44, shifting rows 43:182 down to 45:184.
$ws.Range("A43:A44").EntireRow.Insert()

# Seed the new rows with a copy of the row pair that now lives at 45:46
# (identical to the original 43:44 content) so every column keeps its
# existing formatting/values.
$ws.Range("A45:R46").Copy()
$ws.Range("A43").PasteSpecial()

# Only the date differs for the newly added pair.
$ws.Cells.Item(43, 4).Value = 44525
$ws.Cells.Item(44, 4).Value = 44525
